$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 72
$ws.Range("F6").Value = 836
$ws.Range("F7").Value = 419
$ws.Range("F8").Value = 4693
$ws.Range("F9").Value = 4693
$ws.Range("F11").Value = 121
$ws.Range("F15").Value = 116
$ws.Range("F16").Value = 7478
$ws.Range("F18").Value = 127
$ws.Range("F22").Value = 1359
$ws.Range("F25").Value = 1745
$ws.Range("F27").Value = 2089
$ws.Range("F28").Value = 6171
$ws.Range("F33").Value = 448
$ws.Range("F34").Value = 6407
$ws.Range("F37").Value = 98
$ws.Range("F40").Value = 13
$ws.Range("F41").Value = 2455
$ws.Range("F43").Value = 59
$ws.Range("F44").Value = 1115
$ws.Range("F46").Value = 437
$ws.Range("F47").Value = 2141
$ws.Range("F48").Value = 44
$ws.Range("F49").Value = 1075

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 232
$ws.Range("F14").Value = 23

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1445

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1445
$ws.Range("F5").Value = 232
$ws.Range("F6").Value = 72
$ws.Range("F8").Value = 419
$ws.Range("F9").Value = 4693
$ws.Range("F10").Value = 4693
$ws.Range("F12").Value = 121
$ws.Range("F16").Value = 116
$ws.Range("F17").Value = 7478
$ws.Range("F19").Value = 127
$ws.Range("F21").Value = 1359
$ws.Range("F24").Value = 1745
$ws.Range("F26").Value = 2089
$ws.Range("F29").Value = 6171
$ws.Range("F35").Value = 448
$ws.Range("F36").Value = 6407
$ws.Range("F39").Value = 98
$ws.Range("F42").Value = 2455
$ws.Range("F44").Value = 1115
$ws.Range("F46").Value = 437
$ws.Range("F48").Value = 2141
$ws.Range("F49").Value = 44
$ws.Range("F50").Value = 23
